$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.663.29"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").Value = "3.841.24"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "700.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.27%  "
$ws.Range("D7").Value = "3.840.24"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").Value = "4.476.41"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "3.813.11"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "71.632.76"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "489.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.719"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.14%  "
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").Value = "3.986.00"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +9.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.183"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.42%  "
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("D37").Value = "3.792.32"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000312"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "163.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.35%  "
